# Add season-record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the styling already used by the other header cells (bold, centered,
# top-aligned, thin border) so the new header cells look consistent.
$headerRange = $ws.Range("AC1:AE1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Every player row (2-41) shares the same season record for this team/year.
$lastRow = 41
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 53
    $ws.Cells.Item($r, 30).Value = 61
    $ws.Cells.Item($r, 31).Value = 0
}
